$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 1.6
$ws.Range("I4").Value = "所有 CPS 乘以 1.6。"

$ws.Range("G5").Value = 3
$ws.Range("I5").Value = "点击产量翻至 3 倍，并更快上浮。"

$ws.Range("G7").Value = 3
$ws.Range("I7").Value = "工厂生产效率提升 200%。"

$ws.Range("G8").Value = 50000
$ws.Range("I8").Value = "额外 +50,000 CPS。"
